$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $cell = $ws.Range($cellRef)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue 'D2' '70.469.06'
Set-TextValue 'E2' '  +0.20%  '
Set-TextValue 'D3' '3.609.39'
Set-TextValue 'E3' '  -0.59%  '
Set-TextValue 'E4' '  +0.03%  '
Set-TextValue 'D5' '580.42'
Set-TextValue 'E5' '  -2.30%  '
Set-TextValue 'D6' '189.85'
Set-TextValue 'E6' '  -2.94%  '
Set-TextValue 'D7' '3.605.26'
Set-TextValue 'E7' '  -0.51%  '
Set-TextValue 'E8' '  -2.20%  '
Set-TextValue 'E9' '  +0.09%  '
Set-TextValue 'D10' '0.187'
Set-TextValue 'E10' '  +4.15%  '
Set-TextValue 'D11' '0.661'
Set-TextValue 'E11' '  -1.38%  '
Set-TextValue 'D12' '55.99'
Set-TextValue 'E12' '  -5.00%  '
Set-TextValue 'E13' '  +7.22%  '
Set-TextValue 'D14' '9.69'
Set-TextValue 'E14' '  -2.47%  '
Set-TextValue 'D15' '4.188.95'
Set-TextValue 'E15' '  -0.57%  '
Set-TextValue 'D16' '19.81'
Set-TextValue 'E16' '  -0.13%  '
Set-TextValue 'D17' '3.602.24'
Set-TextValue 'E17' '  -1.01%  '
Set-TextValue 'D18' '70.451.18'
Set-TextValue 'E18' '  +0.22%  '
Set-TextValue 'D19' '12.67'
Set-TextValue 'E19' '  -0.36%  '
Set-TextValue 'E20' '  +0.04%  '
Set-TextValue 'E21' '  -2.15%  '
Set-TextValue 'D22' '488.66'
Set-TextValue 'E22' '  +0.17%  '
Set-TextValue 'D23' '19.38'
Set-TextValue 'E23' '  +0.47%  '
Set-TextValue 'D24' '4.86'
Set-TextValue 'E24' '  -8.70%  '
Set-TextValue 'D25' '96.69'
Set-TextValue 'E25' '  +5.76%  '
Set-TextValue 'D26' '4.36'
Set-TextValue 'E26' '  -2.43%  '
Set-TextValue 'D27' '2.98'
Set-TextValue 'E27' '  -5.89%  '
Set-TextValue 'D28' '11.07'
Set-TextValue 'E28' '  -3.77%  '
Set-TextValue 'D29' '9.40'
Set-TextValue 'E29' '  -2.39%  '
Set-TextValue 'D30' '32.24'
Set-TextValue 'E30' '  -2.23%  '
Set-TextValue 'D31' '7.66'
Set-TextValue 'E31' '  -4.06%  '
Set-TextValue 'D32' '12.19'
Set-TextValue 'E32' '  -0.70%  '
Set-TextValue 'D33' '65.77'
Set-TextValue 'E33' '  -0.27%  '
Set-TextValue 'E34' '  -3.27%  '
Set-TextValue 'D35' '570.61'
Set-TextValue 'E35' '  -8.95%  '
Set-TextValue 'D36' '38.51'
Set-TextValue 'E36' '  -5.80%  '
Set-TextValue 'D37' '0.0₃0812'
Set-TextValue 'E37' '  -1.13%  '
Set-TextValue 'E38' '  +0.15%  '
Set-TextValue 'D39' '0.396'
Set-TextValue 'E39' '  -4.37%  '
Set-TextValue 'D40' '3.30'
Set-TextValue 'E40' '  +13.52%  '
Set-TextValue 'E41' '  +5.87%  '
Set-TextValue 'D42' '3.49'
Set-TextValue 'E42' '  -2.68%  '
Set-TextValue 'D43' '0.138'
Set-TextValue 'E43' '  -6.00%  '
Set-TextValue 'D44' '3.03'
Set-TextValue 'E44' '  -3.94%  '
Set-TextValue 'D45' '3.52'
Set-TextValue 'E45' '  +6.36%  '
Set-TextValue 'D46' '3.224.14'
Set-TextValue 'E46' '  -2.25%  '
Set-TextValue 'E47' '  -1.87%  '
Set-TextValue 'D48' '9.75'
Set-TextValue 'E48' '  +5.78%  '
Set-TextValue 'E49' '  -0.22%  '
Set-TextValue 'B50' 'LidoDAOToken'
Set-TextValue 'C50' 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 'D50' '3.24'
Set-TextValue 'E50' '  -3.46%  '
Set-TextValue 'B51' 'FirstDigitalUSD'
Set-TextValue 'C51' 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 'D51' '0.999'
Set-TextValue 'E51' '  -0.16%  '
